# Add two new hydrogen production pathways ("electrolysis with guaranteed
# clean electricity" and "natural gas reforming with CCS") to the
# HPPECbP sheet, following the pattern of the existing pathway rows
# (label in column A, formula "=$B$2" copied across columns B:AI).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HPPECbP")

# Row 7: electrolysis with guaranteed clean electricity
$ws.Range("A7").Value = "electrolysis with guaranteed clean electricity"
$ws.Range("A7").Font.Bold = $true
$ws.Range("B7:AI7").Formula = "=`$B`$2"

# Row 8: natural gas reforming with CCS
$ws.Range("A8").Value = "natural gas reforming with CCS"
$ws.Range("A8").Font.Bold = $true
$ws.Range("B8:AI8").Formula = "=`$B`$2"

# Reflect the author's on-screen selection/scroll position at save time
# without disturbing which sheet tab is active in the workbook.
$active = $wb.ActiveSheet
$ws.Activate() | Out-Null
$ws.Range("B6:AI8").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 7
$active.Activate() | Out-Null

$wb.Save() | Out-Null
